$d = $word.ActiveDocument

$replacements = @(
    @("2025-05-26 Monday", "2025-05-27 Tuesday"),
    @("21×82=1722", "72×70=5040"),
    @("22×21=462", "55×38=2090"),
    @("82×99=8118", "67×46=3082"),
    @("65×31=2015", "69×63=4347"),
    @("41×71=2911", "64×87=5568"),
    @("11×44=484", "29×51=1479"),
    @("58×48=2784", "53×62=3286"),
    @("52×42=2184", "72×78=5616"),
    @("47×80=3760", "94×74=6956"),
    @("47×15=705", "53×25=1325"),
    @("12×62=744", "42×87=3654"),
    @("96×75=7200", "84×35=2940"),
    @("50×81=4050", "51×97=4947"),
    @("76×57=4332", "71×45=3195"),
    @("98×94=9212", "86×31=2666"),
    @("72×61=4392", "66×38=2508"),
    @("71×94=6674", "69×71=4899"),
    @("46×19=874", "82×93=7626"),
    @("47×29=1363", "23×22=506"),
    @("50×74=3700", "59×86=5074"),
    @("24×28=672", "97×56=5432"),
    @("49×56=2744", "96×68=6528"),
    @("21×66=1386", "26×23=598"),
    @("66×40=2640", "42×39=1638"),
    @("48×34=1632", "95×99=9405")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
